# Fixed Stimulus Absolute Timestamps
# Rename task-order sheets and update their stimulus-file / condition values
# with the refreshed absolute timestamps.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-1650477825022879"
$ws1.Range("B2").Value = "go_stims-16504778249868786.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778250059116.csv"
$ws1.Range("B4").Value = "go_stims-1650477825006879.csv"
$ws1.Range("B5").Value = "GNG_stims-16504778250219145.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16504778280439088"
$ws2.Range("B2").Value = "TB-16504778270649128.csv"
$ws2.Range("B3").Value = "OB-16504778259638784.csv"
$ws2.Range("B4").Value = "TB-16504778280139093.csv"
$ws2.Range("B5").Value = "OB-1650477825883911.csv"
$ws2.Range("B6").Value = "ZB-match_5-1650477825443883.csv"
$ws2.Range("B7").Value = "ZB-match_4-16504778253148754.csv"
$ws2.Range("B8").Value = "TB-16504778267069166.csv"
$ws2.Range("B9").Value = "OB-1650477826121879.csv"
$ws2.Range("B10").Value = "ZB-match_9-16504778251038754.csv"

# --- Sheet 3: RS ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16504778280458803"
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-1650477828091911"
$ws4.Range("B2").Value = "MM_stims-16504778280599139.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778280478787.csv"
$ws4.Range("B4").Value = "MM_stims-16504778280759106.csv"
$ws4.Range("B5").Value = "ZM_stims-1650477828060881.csv"
$ws4.Range("B6").Value = "MM_stims-1650477828091911.csv"
$ws4.Range("B7").Value = "ZM_stims-165047782807688.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-16504778281558774"
$ws5.Range("B2").Value = "vSAT_stims-16504778281399105.csv"
$ws5.Range("B3").Value = "SAT_stims-16504778280958784.csv"
$ws5.Range("B4").Value = "vSAT_stims-16504778281239138.csv"
$ws5.Range("B5").Value = "SAT_stims-16504778281079185.csv"
